$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns.
# D-column price values are forced to stay text (leading apostrophe)
# so Excel doesn't coerce them into floating-point numbers, then the
# cell style is reset to Normal so no stray NumberFormat/quote-prefix
# formatting is left behind on the cell.

$ws.Range("D2").Value = "'59.967.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "'2.417.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'551.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").Value = "'137.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +4.02%  "

$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("D10").Value = "'5.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.70%  "

$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "'25.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "'2.848.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "'59.924.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").Value = "'2.427.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "'11.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("D20").Value = "'328.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("D21").Value = "'6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.06%  "

$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "'65.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("E24").Value = "  +3.54%  "

$ws.Range("D25").Value = "'8.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").Value = "'0.0₃0776"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("E29").Value = "  -3.03%  "

$ws.Range("D30").Value = "'168.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("E31").Value = "  -4.45%  "

$ws.Range("D32").Value = "'18.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -0.99%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -2.56%  "

$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("D39").Value = "'320.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").Value = "'0.405"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("D42").Value = "'140.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.66%  "

$ws.Range("E43").Value = "  +0.57%  "

$ws.Range("D44").Value = "'19.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "

$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("D46").Value = "'0.578"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("D48").Value = "'0.387"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.96%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  -3.99%  "

$ws.Range("E51").Value = "  -1.10%  "
